# Daily attendance processing - 2026-01-04 20:36:41
#
# The "Recorded By" column (G) lists the users who touched a session
# record, separated by ", ". This pass normalizes the ordering so the
# real recorder's email comes before the "System" marker (and keeps any
# trailing lowercase "system" marker last), e.g.
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"           -> "backup@backdoor.com, System"
#   "System, backup@backdoor.com, system"   -> "backup@backdoor.com, System, system"
#   "admin@admin.com, dnasr281@gmail.com"   -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ", "
    if ($parts.Count -ge 2) {
        $first = $parts[0]
        $second = $parts[1]

        # "System" should always yield the lead slot to the other
        # identity unless that identity is "admin@admin.com" (System
        # keeps priority over the generic admin account).
        $swap = $false
        if ($first -eq "System" -and $second -ne "admin@admin.com") {
            $swap = $true
        } elseif ($first -eq "admin@admin.com" -and $second -ne "System") {
            $swap = $true
        }

        if ($swap) {
            $parts[0] = $second
            $parts[1] = $first
            $cell.Value2 = ($parts -join ", ")
        }
    }
}
